$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to stay plain text, then drop back to the default (unstyled) style
# so no stray formatting is introduced.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '69.078.40'
$ws.Range('E2').Value = '  -3.09%  '
$ws.Range('D3').Value = '3.516.52'
$ws.Range('E3').Value = '  -4.90%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue 'D5' '577.83'
$ws.Range('E5').Value = '  -0.89%  '
Set-TextValue 'D6' '171.11'
$ws.Range('E6').Value = '  -3.74%  '
$ws.Range('D7').Value = '3.508.77'
$ws.Range('E7').Value = '  -4.77%  '
Set-TextValue 'D8' '0.608'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  -5.69%  '
Set-TextValue 'D11' '6.75'
$ws.Range('E11').Value = '  -2.25%  '
Set-TextValue 'D12' '0.582'
$ws.Range('E12').Value = '  -4.38%  '
Set-TextValue 'D13' '46.90'
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('E14').Value = '  -4.64%  '
$ws.Range('D15').Value = '4.081.76'
$ws.Range('E15').Value = '  -4.92%  '
$ws.Range('E16').Value = '  -5.35%  '
Set-TextValue 'D17' '622.78'
$ws.Range('E17').Value = '  -7.80%  '
$ws.Range('D18').Value = '69.052.57'
$ws.Range('E18').Value = '  -3.21%  '
$ws.Range('D19').Value = '3.484.43'
$ws.Range('E19').Value = '  -5.82%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -3.03%  '
$ws.Range('E22').Value = '  -3.81%  '
$ws.Range('E23').Value = '  -6.17%  '
Set-TextValue 'D24' '15.88'
$ws.Range('E24').Value = '  -8.83%  '
Set-TextValue 'D25' '97.33'
$ws.Range('E25').Value = '  -4.53%  '
Set-TextValue 'D26' '3.79'
$ws.Range('E26').Value = '  -4.50%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  -6.74%  '
Set-TextValue 'D29' '9.32'
$ws.Range('E29').Value = '  -9.48%  '
Set-TextValue 'D30' '32.58'
$ws.Range('E30').Value = '  -7.34%  '
Set-TextValue 'D31' '3.15'
$ws.Range('E31').Value = '  -7.95%  '
Set-TextValue 'D32' '8.53'
$ws.Range('E32').Value = '  -7.07%  '
$ws.Range('E33').Value = '  -7.13%  '
Set-TextValue 'D34' '6.98'
$ws.Range('E34').Value = '  -6.23%  '
Set-TextValue 'D35' '634.00'
$ws.Range('E35').Value = '  +8.05%  '
Set-TextValue 'D36' '10.73'
$ws.Range('E36').Value = '  -3.98%  '
$ws.Range('E37').Value = '  -5.56%  '
Set-TextValue 'D38' '56.73'
Set-TextValue 'D39' '3.41'
$ws.Range('E39').Value = '  -16.34%  '
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('E41').Value = '  -2.10%  '
Set-TextValue 'D42' '0.136'
$ws.Range('E42').Value = '  -6.00%  '
$ws.Range('D43').Value = '3.379.03'
$ws.Range('E43').Value = '  -8.28%  '
Set-TextValue 'D44' '0.326'
$ws.Range('E44').Value = '  -6.71%  '
Set-TextValue 'D45' '32.83'
$ws.Range('E45').Value = '  -7.70%  '
$ws.Range('D46').Value = '0.0₃0687'
$ws.Range('E46').Value = '  -10.12%  '
$ws.Range('E47').Value = '  -7.43%  '
$ws.Range('E48').Value = '  -4.76%  '
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('E50').Value = '  +14.74%  '
Set-TextValue 'D51' '132.16'
